# "Keywords - collect column"
#
# Adds a new "collect" column (B) next to the existing keywords column (A),
# marking every keyword row with "s", and inserts a new keyword
# ("Rotulo alimento") just before the last existing row ("Embalagem alimento").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + flag values for column B, written before the row insert so the
# shared-string table picks up "collect" / "s" ahead of the new keyword text
# (matches the string order produced by the original edit).
$ws.Range("B1").Value = "collect"
$ws.Range("B2:B9").Value = "s"

# Insert a new row 9 (pushes the old "Embalagem alimento" row down to row 10)
# and populate it with the new keyword plus its "collect" flag.
$ws.Rows(9).Insert()
$ws.Range("B9").Value = "s"
$ws.Range("A9").Value = "Rotulo alimento"

# Leave the selection on A2, matching the saved view state.
[void]$ws.Range("A2").Select()
